# Day 25 Code examples edits
#
# 1) Slide 15, shape 4 ("Java Dates ..." text box): expand the explanation of
#    what the long value passed to the Date constructor represents, fix the
#    "Thes" typo -> "This", and correct the epoch year 1990 -> 1970.
# 2) Slide 2, shape 2 (Maven bullet list): change "Java projects" to
#    "programming projects" in the second bullet.

$p = $ppt.ActivePresentation

# --- Edit 1: slide 15 ---------------------------------------------------
$slide15 = $p.Slides.Item(15)
$shape15 = $slide15.Shapes.Item(4)
$tr15 = $shape15.TextFrame.TextRange
$para1 = $tr15.Paragraphs(1)

# Rewrite the first part of the sentence (covers the original "Java Dates
# are created using a constructor which takes a long value. Thes" run plus
# the leading " value represents milliseconds since Jan 01, 19" of the third
# run) so it reads "...a long value. This value represents milliseconds
# since Jan 01" - this also fixes the "Thes" -> "This" typo.
$run1 = $para1.Characters(1, 120)
$run1.Text = "Java Dates are created using a constructor which takes a long value. This value represents milliseconds since Jan 01"

# The remaining text right after that is now "90 (based on the United ...".
# Turn the leftover "90 " into the new ", 1970 " (fixing 1990 -> 1970).
$run2 = $para1.Characters(117, 3)
$run2.Text = ", 1970 "

# --- Edit 2: slide 2 ------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$mavenPara = $tr2.Paragraphs(2)

$mavenRun = $mavenPara.Characters(1, $mavenPara.Length)
$mavenRun.Text = "Maven is a tool that can be used for building and managing programming projects"
